# Added a check for multiple bot replies to same comment.
# The "opinions" source list (used by the bot when picking a random
# reply) gets a few topics renamed/recapitalized and a batch of new
# drag-queen / food / fandom topics appended.
#
# Column A is a flat list of opinion topics, one per row. We only ever
# touch Range.Value so the engine manages the shared-string table itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the brand-new topics first (rows 40-46), in the same order
#     they appear in the final sheet, so the shared-string table grows in
#     the right sequence before the renames below reuse/retire old slots.
$ws.Range("A40").Value = "RuPaul"
$ws.Range("A41").Value = "Valentina"
$ws.Range("A42").Value = "Willam"
$ws.Range("A43").Value = "Sharon Needles"
$ws.Range("A44").Value = "Alyssa Edwards"
$ws.Range("A45").Value = "Whataburger"
$ws.Range("A46").Value = "In-N-Out"

# --- Re-casing / renaming of existing topics (same row position, new text).
#     "Football" -> "football", "Spongebob" -> "SpongeBob", "Poop" -> "poop".
$ws.Range("A33").Value = "football"
$ws.Range("A15").Value = "SpongeBob"
$ws.Range("A12").Value = "poop"

# --- Two more brand-new topics appended at the very end.
$ws.Range("A47").Value = "Star Wars"
$ws.Range("A48").Value = "Star Trek"

# --- Update the view state to match: selection sitting one row past the
#     new last row, scrolled down so row 13 is at the top of the window.
$ws.Activate() | Out-Null
$ws.Range("A49").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
